# Actualización automatica mar abr  6 17:33:39 CEST 2021
#
# Two rows were removed from the mapping table:
#   - "Sudan del Sur" (row 51)
#   - "Bahamas" (the last data row)
# Deleting whole rows shifts everything below them up by one, which is
# exactly what the target workbook shows (e.g. what used to be row 52
# "Benin" becomes the new row 51, etc.), and the two now-unused trailing
# rows disappear entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "Sudan del Sur" - row 51.
$ws.Rows.Item(51).Delete()

# After the shift above, the former last row ("Bahamas") is now row 141.
$ws.Rows.Item(141).Delete()
